$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 43, shifting existing rows 43:59 down to 44:60.
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row 43 with the new record.
$ws.Cells.Item(43, 1).Value = 5
$ws.Cells.Item(43, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(43, 3).Value = "Maule"
$ws.Cells.Item(43, 4).Value = 44511
$ws.Cells.Item(43, 4).NumberFormat = $ws.Cells.Item(44, 4).NumberFormat
$ws.Cells.Item(43, 5).Value = 7
$ws.Cells.Item(43, 6).Value = 100112022
$ws.Cells.Item(43, 7).Value = "Arveja Verde"
$ws.Cells.Item(43, 8).Value = "Sin especificar"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 500
$ws.Cells.Item(43, 11).Value = 14000
$ws.Cells.Item(43, 12).Value = 14000
$ws.Cells.Item(43, 13).Value = 14000
$ws.Cells.Item(43, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(43, 15).Value = "Región del Maule"
$ws.Cells.Item(43, 16).Value = 560
$ws.Cells.Item(43, 17).Value = 25
$ws.Cells.Item(43, 18).Value = "Hortaliza"
